# Update metricas_recorrencia_trimestral data for 2025Q3 (row 21)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 115
$ws.Range("D21").Value = 106
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = 30.37249283667622
